# Auto-generated PowerShell COM-interop script
# Applies the 2025-12-12 01:24:35 append/update to the lancers sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks so we can rebuild them cleanly without duplicates
$ws.Hyperlinks.Delete()

$timestamp = '2025-12-12 01:24:35'

# Row 2
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = '【急募】AIチャットボット開発のプロフェッショナルを探しています!'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5451734'
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5451734')

# Row 3
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Range("G3").Value = 178
$ws.Range("H3").Value = '★bot ◆ツール'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5405023')

# Row 4
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = 'Webシステム開発'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5451859'
$ws.Range("G4").Value = 118
$ws.Range("H4").Value = '◆開発,システム開発'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5451859')

# Row 5
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = '【Flutterエンジニア募集】Androidアプリ開発のパートナーを探しています'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5452211'
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = '◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5452211')

# Row 6
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = '自動出品システムの開発'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5451514'
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5451514')

# Row 7
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = '【愛知県近辺 在住の方希望 / リモート相談可能】経験豊富なWebフロントエンド開発エンジニア募集!'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5451972'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5451972')

# Row 8
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = '【フルリモート】WordPressサイトの構築・運用'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5449760'
$ws.Range("G8").Value = 58
$ws.Range("H8").Value = '◇サイト ○WordPress'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5449760')

# Row 9
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = 'Amazonの購入アカウントから必要な情報のスクレイピング→スプレッドシートに記入をしたい。'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5452210'
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = '◆スクレイピング'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5452210')

# Row 10
$ws.Range("A10").Value = $timestamp
$ws.Range("B10").Value = '【R/Shiny】高齢者評価アプリ 機能追加・UI改修依頼'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5452159'
$ws.Range("G10").Value = 38
$ws.Range("H10").Value = '◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5452159')

# Row 11
$ws.Range("A11").Value = $timestamp
$ws.Range("B11").Value = 'Access DB家賃管理SYSを最新Access で稼働できるように'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5451626'
$ws.Range("G11").Value = 38
$ws.Range("H11").Value = '◇管理'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5451626')

# Row 12
$ws.Range("A12").Value = $timestamp
$ws.Range("B12").Value = '【小規模・短納期・急募】アプリMatrixifyを用いたデータ移行検証・マッピング担当募集'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5451926'
$ws.Range("G12").Value = 33
$ws.Range("H12").Value = '◇アプリ'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5451926')

# Row 13
$ws.Range("A13").Value = $timestamp
$ws.Range("B13").Value = '【急募】フルスクラッチECサイトのグーグルサーチコンソール設定依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5452161'
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = '◇サイト'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5452161')

# Row 14
$ws.Range("A14").Value = $timestamp
$ws.Range("B14").Value = '進行管理およびチームディレクションを担当'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '~ 5,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = '◇管理'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5418064')

# Row 15
$ws.Range("A15").Value = $timestamp
$ws.Range("B15").Value = 'Rubyの暗号化機能のPHP化'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5451714'
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = '○PHP'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5451714')

# Row 16
$ws.Range("A16").Value = $timestamp
$ws.Range("B16").Value = '【オンライン講師募集】バックエンドの基礎を教えていただける方'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5451420'
$ws.Range("G16").Value = 18
$ws.Range("H16").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5451420')

# Row 17
$ws.Range("A17").Value = $timestamp
$ws.Range("B17").Value = '注目 限定公開 PR 限定公開の仕事'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5450323'
$ws.Range("G17").Value = 13
$ws.Range("H17").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5450323')

# Row 18
$ws.Range("A18").Value = $timestamp
$ws.Range("B18").Value = '【急募】desknetスタンダード版からNeo版への移行サポート'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5451838'
$ws.Range("G18").Value = 13
$ws.Range("H18").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5451838')

# Row 19
$ws.Range("A19").Value = $timestamp
$ws.Range("B19").Value = 'Xの運用代行'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5451931'
$ws.Range("G19").Value = 10
$ws.Range("H19").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5451931')
